# Updates cryptos list data (prices / 1h volume change) to refreshed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.508.99"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "3.907.73"
$ws.Range("E3").Value = "  +3.64%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "602.37"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").Value = "164.28"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "3.904.94"
$ws.Range("E7").Value = "  +3.61%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -2.59%  "
$ws.Range("E10").Value = "  -3.95%  "
$ws.Range("D11").Value = "6.36"
$ws.Range("E11").Value = "  +0.16%  "
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "36.72"
$ws.Range("E14").Value = "  -2.52%  "
$ws.Range("D15").Value = "4.558.74"
$ws.Range("E15").Value = "  +3.57%  "
$ws.Range("D16").Value = "3.930.41"
$ws.Range("E16").Value = "  +3.86%  "
$ws.Range("D17").Value = "68.659.83"
$ws.Range("E17").Value = "  -0.88%  "
$ws.Range("D18").Value = "7.39"
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").Value = "16.98"
$ws.Range("E20").Value = "  -3.73%  "
$ws.Range("D21").Value = "11.22"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("D22").Value = "483.32"
$ws.Range("E22").Value = "  -1.75%  "
$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D23").Value = "0.0000169"
$ws.Range("E23").Value = "  +12.18%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.716"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").Value = "84.29"
$ws.Range("E25").Value = "  -0.63%  "
$ws.Range("E26").Value = "  -1.69%  "
$ws.Range("D27").Value = "11.95"
$ws.Range("E27").Value = "  -2.65%  "
$ws.Range("D28").Value = "10.06"
$ws.Range("E28").Value = "  -0.41%  "
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("E30").Value = "  -1.59%  "
$ws.Range("D31").Value = "4.055.80"
$ws.Range("E31").Value = "  +3.62%  "
$ws.Range("D32").Value = "7.82"
$ws.Range("E32").Value = "  -3.84%  "
$ws.Range("E33").Value = "  -2.76%  "
$ws.Range("D34").Value = "31.87"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "3.849.23"
$ws.Range("E35").Value = "  +3.37%  "
$ws.Range("E36").Value = "  -1.27%  "
$ws.Range("D37").Value = "1.04"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "5.86"
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -2.80%  "
$ws.Range("D43").Value = "431.09"
$ws.Range("E43").Value = "  +1.13%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("D45").Value = "1.97"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("D47").Value = "'8.40"
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "26.41"
$ws.Range("E48").Value = "  +11.21%  "
$ws.Range("D49").Value = "141.83"
$ws.Range("E49").Value = "  -0.23%  "
$ws.Range("D50").Value = "2.812.25"
$ws.Range("E50").Value = "  +0.19%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "0.0351"
$ws.Range("E51").Value = "  -0.38%  "
